$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells.
# For Price (column D) cells whose new value would parse as a plain number
# (single decimal point, digits only), force the cell to Text format first
# so Excel stores it as a text string (matching the sheet-wide inline-string
# convention) instead of silently coercing it to a numeric value.

$ws.Range("D2").Value = '62.467.66'
$ws.Range("E2").Value = '  +2.81%  '

$ws.Range("D3").Value = '2.428.48'
$ws.Range("E3").Value = '  +3.64%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.08'
$ws.Range("E5").Value = '  +2.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.83'
$ws.Range("E6").Value = '  +5.65%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  +2.03%  '

$ws.Range("D9").Value = '2.429.77'
$ws.Range("E9").Value = '  +3.77%  '

$ws.Range("E10").Value = '  +5.64%  '

$ws.Range("E11").Value = '  -0.37%  '

$ws.Range("E12").Value = '  +2.28%  '

$ws.Range("E13").Value = '  +3.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.35'
$ws.Range("E14").Value = '  +7.30%  '

$ws.Range("E15").Value = '  +10.25%  '

$ws.Range("D16").Value = '2.866.66'
$ws.Range("E16").Value = '  +3.61%  '

$ws.Range("D17").Value = '62.318.06'
$ws.Range("E17").Value = '  +2.65%  '

$ws.Range("D18").Value = '2.427.85'
$ws.Range("E18").Value = '  +3.70%  '

$ws.Range("E19").Value = '  +4.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '324.88'
$ws.Range("E20").Value = '  +2.12%  '

$ws.Range("E21").Value = '  +1.64%  '

$ws.Range("E22").Value = '  +3.60%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("E24").Value = '  +6.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.08'
$ws.Range("E25").Value = '  +2.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.08'
$ws.Range("E26").Value = '  +7.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '569.72'
$ws.Range("E27").Value = '  +14.75%  '

$ws.Range("D28").Value = '2.548.34'
$ws.Range("E28").Value = '  +3.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("E30").Value = '  +10.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.44'
$ws.Range("E31").Value = '  +6.45%  '

$ws.Range("E32").Value = '  +6.84%  '

$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").Value = '  +4.55%  '

$ws.Range("E35").Value = '  +5.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.72'
$ws.Range("E36").Value = '  +9.36%  '

$ws.Range("E37").Value = '  +6.12%  '

$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("E39").Value = '  +2.81%  '

$ws.Range("E40").Value = '  +4.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.81'
$ws.Range("E41").Value = '  +1.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '150.67'
$ws.Range("E42").Value = '  +5.66%  '

$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("E44").Value = '  +2.93%  '

$ws.Range("E45").Value = '  +16.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '151.47'
$ws.Range("E46").Value = '  +6.44%  '

$ws.Range("E47").Value = '  +3.17%  '

$ws.Range("E48").Value = '  +5.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.46'
$ws.Range("E49").Value = '  +7.80%  '

$ws.Range("E50").Value = '  +4.48%  '

$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0918'
$ws.Range("E51").Value = '  +2.16%  '
